# Update column G (K = strikeouts) with newly regenerated simulated values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 1
    3 = 0
    4 = 0
    5 = 1
    6 = 1
    7 = 3
    8 = 3
    9 = 1
    10 = 0
    11 = 3
    12 = 1
    13 = 1
    14 = 0
    15 = 3
    16 = 0
    17 = 1
    18 = 0
    19 = 2
    20 = 3
    21 = 1
    22 = 0
    23 = 0
    24 = 2
    25 = 3
    26 = 1
    27 = 1
    28 = 0
    29 = 2
    30 = 1
    31 = 2
    32 = 2
    33 = 0
    34 = 1
    35 = 0
    36 = 0
    37 = 1
    38 = 1
    39 = 3
    40 = 2
    41 = 2
    42 = 1
    43 = 1
    44 = 2
    45 = 1
    46 = 1
    48 = 1
    49 = 1
    50 = 0
    51 = 0
    52 = 2
    53 = 0
    54 = 0
    55 = 1
    56 = 1
    57 = 2
    59 = 2
    60 = 1
    61 = 1
    62 = 2
    63 = 2
    64 = 1
    65 = 0
    66 = 1
    67 = 1
    68 = 1
    69 = 1
    70 = 1
    71 = 0
    72 = 2
    73 = 1
    74 = 1
    75 = 2
    76 = 1
    77 = 1
    78 = 2
    79 = 2
    80 = 2
    83 = 2
    84 = 3
    85 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
